# Updated cryptos list - refresh Price (D) and Volume(1h) (E) columns
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "29.394.58"
$ws.Range("E2").Value = "  +0.26%  "
$ws.Range("D3").Value = "1.842.01"
$ws.Range("E3").Value = "  -0.11%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.9996"
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "239.37"
$ws.Range("E5").Value = "  -0.31%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.6259"
$ws.Range("E6").Value = "  -0.21%  "
$ws.Range("E7").Value = "  +0.19%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.07438"
$ws.Range("E8").Value = "  -0.68%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "24.99"
$ws.Range("E9").Value = "  +2.36%  "
$ws.Range("E10").Value = "  -0.13%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07720"
$ws.Range("E11").Value = "  -0.17%  "
$ws.Range("D12").Value = "1.840.83"
$ws.Range("E12").Value = "  -0.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.979"
$ws.Range("E13").Value = "  -0.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6759"
$ws.Range("E14").Value = "  -0.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.00001028"
$ws.Range("E15").Value = "  -2.09%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "81.89"
$ws.Range("E16").Value = "  -0.17%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.229"
$ws.Range("E17").Value = "  +0.80%  "
$ws.Range("D18").Value = "29.397.61"
$ws.Range("E18").Value = "  +0.05%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "233.17"
$ws.Range("E19").Value = "  +1.68%  "
$ws.Range("E20").Value = "  +0.13%  "
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.332"
$ws.Range("E22").Value = "  -2.07%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "158.18"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "8.487"
$ws.Range("E25").Value = "  +0.85%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1352"
$ws.Range("E26").Value = "  -1.62%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "17.33"
$ws.Range("E27").Value = "  -1.09%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.07181"
$ws.Range("E28").Value = "  +12.17%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.468"
$ws.Range("E29").Value = "  +3.11%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.483"
$ws.Range("E30").Value = "  +0.06%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "4.046"
$ws.Range("E31").Value = "  -1.16%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.030"
$ws.Range("E32").Value = "  -1.42%  "
$ws.Range("E33").Value = "  -0.61%  "
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.6978"
$ws.Range("E35").Value = "  +0.03%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.576"
$ws.Range("E36").Value = "  -0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.01840"
$ws.Range("E37").Value = "  +0.40%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.934"
$ws.Range("E38").Value = "  +3.36%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.816"
$ws.Range("E39").Value = "  -0.79%  "
$ws.Range("D40").Value = "1.234.68"
$ws.Range("E40").Value = "  -2.76%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.9572"
$ws.Range("E41").Value = "  +4.74%  "
$ws.Range("E42").Value = "  +0.24%  "
$ws.Range("D43").Value = "2.011.88"
$ws.Range("E43").Value = "  +0.16%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "101.01"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "65.50"
$ws.Range("E45").Value = "  -1.11%  "
$ws.Range("E46").Value = "  +4.58%  "
$ws.Range("E47").Value = "  -0.28%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "6.968"
$ws.Range("E48").Value = "  -1.46%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.885"
$ws.Range("E49").Value = "  -1.45%  "
$ws.Range("E50").Value = "  -2.35%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3902"
$ws.Range("E51").Value = "  -1.50%  "
